$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 70

# Column A holds a plain date-like string (e.g. "11/09/2025" in the row
# above), not a real Excel date. Force text interpretation so it isn't
# auto-converted to a date serial, then restore the default ("Normal")
# style so the new row matches the unstyled look of the other data rows.
$cellA = $ws.Range("A" + $row)
$cellA.NumberFormat = "@"
$cellA.Value = "11/10/2025"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.1940694112342378
$ws.Cells.Item($row, 3).Value = 0.8059305887657622
